$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dbExcel query in B2 is rewritten to alias its RETURN columns with
# human-readable (backtick-quoted where they contain spaces) names instead
# of the raw snake_case property names.
$newQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN ['Clinical Trial','Transcriptomics'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

$ws.Cells.Item(2, 2).Value = $newQuery

# The cell keeps wrapping, and with the longer text the row now needs to be
# taller to show it in full.
$ws.Rows.Item(2).RowHeight = 188.5

# Selection/top-left moved from the dbExcel column (whole column C selected)
# to a single-cell selection on the query cell (B2).
$ws.Range("B2").Select() | Out-Null
